$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-10-15 Tuesday" "2024-10-16 Wednesday"

Replace-Text "69×80=" "57×48="
Replace-Text "30×68=" "12×36="
Replace-Text "46×61=" "26×71="
Replace-Text "34×90=" "50×58="
Replace-Text "52×82=" "39×73="

Replace-Text "24×20=" "32×74="
Replace-Text "44×45=" "68×76="
Replace-Text "36×41=" "99×75="
Replace-Text "85×90=" "46×36="
Replace-Text "12×79=" "92×97="

Replace-Text "12×90=" "15×22="
Replace-Text "82×42=" "13×75="
Replace-Text "66×36=" "51×89="
Replace-Text "65×56=" "57×17="
Replace-Text "24×25=" "44×24="

Replace-Text "70×90=" "16×53="
Replace-Text "73×72=" "54×62="
Replace-Text "54×19=" "68×12="
Replace-Text "40×75=" "90×13="
Replace-Text "21×22=" "28×68="

Replace-Text "70×47=" "82×77="
Replace-Text "93×35=" "98×29="
Replace-Text "51×76=" "49×94="
Replace-Text "11×68=" "49×43="
Replace-Text "24×43=" "76×81="
